# fastexcel-test fill/byName.xlsx — "refactor: Update the project's fill
# test demo." Sheet1 holds the fill-by-name template: row 1 is the set of
# header/placeholder labels, row 2 is the template text that gets expanded.
# The Chinese labels/templates are swapped for English ones and the
# corresponding header/empty-marker cells pick up the "宋体" font used
# elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new cell text in the same order the original authoring tool
# appended them to the shared-string table (A1, B1, C2, D2, E1, E2, D1, C1)
# so new unique strings land at the same indices as the target workbook.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C2").Value = "{name} is {number} years old this year"
$ws.Range("D2").Value = "\{name\} ignored, {name}"
$ws.Range("E1").Value = "Empty"
$ws.Range("E2").Value = "Empty{.empty}"
$ws.Range("D1").Value = "Ignored"
$ws.Range("C1").Value = "Complex"

# Header row + the two template cells that used to carry the default font
# now use the 宋体 font (matching the small-print font already used
# elsewhere in this style sheet).
$ws.Range("A1:D1").Font.Name = "宋体"
$ws.Range("C2:D2").Font.Name = "宋体"
$ws.Range("E1:E2").Font.Name = "宋体"

# Restore the cursor to the cell the author left selected.
$ws.Range("C2").Select() | Out-Null
